# Remove negative cross sections (rows whose "value" column I is negative)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 31, 15, 13 (descending order so earlier deletions don't
# shift the row numbers of rows still to be deleted)
$ws.Rows.Item(31).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(13).Delete()
